$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("D3").Value = 44203
$ws.Range("J3").Value = 27
$ws.Range("K3").Value = 7000
$ws.Range("M3").Value = 7556
$ws.Range("P3").Value = 756

# Row 5 updates
$ws.Range("D5").Value = 44775
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 8000
$ws.Range("M5").Value = 8000
$ws.Range("P5").Value = 800
